$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.36733066666667
$ws.Range("H2").Value = 31.101992
$ws.Range("I2").Value = 0.1169328841728879
$ws.Range("J2").Value = 0.1169328841728879
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 98.946724
$ws.Range("N2").Value = 296.840172
$ws.Range("O2").Value = 0.2098009692989996
$ws.Range("P2").Value = 0.2098009692989996
$ws.Range("Q2").Value = 1025.813406091403
$ws.Range("R2").Value = 9232.320654822623
$ws.Range("S2").Value = 0.02453263244239953
$ws.Range("T2").Value = 0.02453263244239954

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.36733066666667
$ws.Range("H3").Value = 31.101992
$ws.Range("I3").Value = 0.1169328841728879
$ws.Range("J3").Value = 0.1169328841728879
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 163.0062356666667
$ws.Range("N3").Value = 489.018707
$ws.Range("O3").Value = 0.345629090707923
$ws.Range("P3").Value = 0.3456290907079231
$ws.Range("Q3").Value = 1689.939545884927
$ws.Range("R3").Value = 15209.45591296434
$ws.Range("S3").Value = 0.04041540643053013
$ws.Range("T3").Value = 0.04041540643053015

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.36733066666667
$ws.Range("H4").Value = 31.101992
$ws.Range("I4").Value = 0.1169328841728879
$ws.Range("J4").Value = 0.1169328841728879
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 65.39610666666668
$ws.Range("N4").Value = 196.18832
$ws.Range("O4").Value = 0.1386621609326595
$ws.Range("P4").Value = 0.1386621609326595
$ws.Range("Q4").Value = 677.9830621259379
$ws.Range("R4").Value = 6101.84755913344
$ws.Range("S4").Value = 0.01621416640350101
$ws.Range("T4").Value = 0.01621416640350102

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 10.36733066666667
$ws.Range("H5").Value = 31.101992
$ws.Range("I5").Value = 0.1169328841728879
$ws.Range("J5").Value = 0.1169328841728879
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 144.2727966666667
$ws.Range("N5").Value = 432.81839
$ws.Range("O5").Value = 0.3059077790604178
$ws.Range("P5").Value = 0.3059077790604179
$ws.Range("Q5").Value = 1495.723789248098
$ws.Range("R5").Value = 13461.51410323288
$ws.Range("S5").Value = 0.03577067889645722
$ws.Range("T5").Value = 0.03577067889645724

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 37.91490933333333
$ws.Range("H6").Value = 113.744728
$ws.Range("I6").Value = 0.4276413904453658
$ws.Range("J6").Value = 0.4276413904453659
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 98.946724
$ws.Range("N6").Value = 296.840172
$ws.Range("O6").Value = 0.2098009692989996
$ws.Range("P6").Value = 0.2098009692989996
$ws.Range("Q6").Value = 3751.556069290357
$ws.Range("R6").Value = 33764.00462361322
$ws.Range("S6").Value = 0.08971957822780967
$ws.Range("T6").Value = 0.08971957822780971

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 37.91490933333333
$ws.Range("H7").Value = 113.744728
$ws.Range("I7").Value = 0.4276413904453658
$ws.Range("J7").Value = 0.4276413904453659
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 163.0062356666667
$ws.Range("N7").Value = 489.018707
$ws.Range("O7").Value = 0.345629090707923
$ws.Range("P7").Value = 0.3456290907079231
$ws.Range("Q7").Value = 6180.366646069633
$ws.Range("R7").Value = 55623.2998146267
$ws.Range("S7").Value = 0.1478053049287036
$ws.Range("T7").Value = 0.1478053049287037

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 37.91490933333333
$ws.Range("H8").Value = 113.744728
$ws.Range("I8").Value = 0.4276413904453658
$ws.Range("J8").Value = 0.4276413904453659
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 65.39610666666668
$ws.Range("N8").Value = 196.18832
$ws.Range("O8").Value = 0.1386621609326595
$ws.Range("P8").Value = 0.1386621609326595
$ws.Range("Q8").Value = 2479.487455019663
$ws.Range("R8").Value = 22315.38709517696
$ws.Range("S8").Value = 0.05929767930340157
$ws.Range("T8").Value = 0.05929767930340159

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 37.91490933333333
$ws.Range("H9").Value = 113.744728
$ws.Range("I9").Value = 0.4276413904453658
$ws.Range("J9").Value = 0.4276413904453659
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 144.2727966666667
$ws.Range("N9").Value = 432.81839
$ws.Range("O9").Value = 0.3059077790604178
$ws.Range("P9").Value = 0.3059077790604179
$ws.Range("Q9").Value = 5470.090004883103
$ws.Range("R9").Value = 49230.81004394792
$ws.Range("S9").Value = 0.1308188279854508
$ws.Range("T9").Value = 0.1308188279854509

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 26.72147866666667
$ws.Range("H10").Value = 80.164436
$ws.Range("I10").Value = 0.3013909433702152
$ws.Range("J10").Value = 0.3013909433702153
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 98.946724
$ws.Range("N10").Value = 296.840172
$ws.Range("O10").Value = 0.2098009692989996
$ws.Range("P10").Value = 0.2098009692989996
$ws.Range("Q10").Value = 2644.002774502555
$ws.Range("R10").Value = 23796.02497052299
$ws.Range("S10").Value = 0.06323211205701104
$ws.Range("T10").Value = 0.06323211205701107

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 26.72147866666667
$ws.Range("H11").Value = 80.164436
$ws.Range("I11").Value = 0.3013909433702152
$ws.Range("J11").Value = 0.3013909433702153
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 163.0062356666667
$ws.Range("N11").Value = 489.018707
$ws.Range("O11").Value = 0.345629090707923
$ws.Range("P11").Value = 0.3456290907079231
$ws.Range("Q11").Value = 4355.767648900473
$ws.Range("R11").Value = 39201.90884010425
$ws.Range("S11").Value = 0.1041694777046506
$ws.Range("T11").Value = 0.1041694777046506

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 26.72147866666667
$ws.Range("H12").Value = 80.164436
$ws.Range("I12").Value = 0.3013909433702152
$ws.Range("J12").Value = 0.3013909433702153
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 65.39610666666668
$ws.Range("N12").Value = 196.18832
$ws.Range("O12").Value = 0.1386621609326595
$ws.Range("P12").Value = 0.1386621609326595
$ws.Range("Q12").Value = 1747.480669176392
$ws.Range("R12").Value = 15727.32602258752
$ws.Range("S12").Value = 0.04179151949324683
$ws.Range("T12").Value = 0.04179151949324685

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 26.72147866666667
$ws.Range("H13").Value = 80.164436
$ws.Range("I13").Value = 0.3013909433702152
$ws.Range("J13").Value = 0.3013909433702153
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 144.2727966666667
$ws.Range("N13").Value = 432.81839
$ws.Range("O13").Value = 0.3059077790604178
$ws.Range("P13").Value = 0.3059077790604179
$ws.Range("Q13").Value = 3855.182458308671
$ws.Range("R13").Value = 34696.64212477804
$ws.Range("S13").Value = 0.0921978341153067
$ws.Range("T13").Value = 0.09219783411530673

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 13.65680433333333
$ws.Range("H14").Value = 40.970413
$ws.Range("I14").Value = 0.154034782011531
$ws.Range("J14").Value = 0.154034782011531
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 98.946724
$ws.Range("N14").Value = 296.840172
$ws.Range("O14").Value = 0.2098009692989996
$ws.Range("P14").Value = 0.2098009692989996
$ws.Range("Q14").Value = 1351.296049092337
$ws.Range("R14").Value = 12161.66444183104
$ws.Range("S14").Value = 0.03231664657177931
$ws.Range("T14").Value = 0.03231664657177932

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 13.65680433333333
$ws.Range("H15").Value = 40.970413
$ws.Range("I15").Value = 0.154034782011531
$ws.Range("J15").Value = 0.154034782011531
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 163.0062356666667
$ws.Range("N15").Value = 489.018707
$ws.Range("O15").Value = 0.345629090707923
$ws.Range("P15").Value = 0.3456290907079231
$ws.Range("Q15").Value = 2226.144265612888
$ws.Range("R15").Value = 20035.29839051599
$ws.Range("S15").Value = 0.0532389016440386
$ws.Range("T15").Value = 0.05323890164403862

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 13.65680433333333
$ws.Range("H16").Value = 40.970413
$ws.Range("I16").Value = 0.154034782011531
$ws.Range("J16").Value = 0.154034782011531
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 65.39610666666668
$ws.Range("N16").Value = 196.18832
$ws.Range("O16").Value = 0.1386621609326595
$ws.Range("P16").Value = 0.1386621609326595
$ws.Range("Q16").Value = 893.1018329084625
$ws.Range("R16").Value = 8037.916496176162
$ws.Range("S16").Value = 0.02135879573251003
$ws.Range("T16").Value = 0.02135879573251004

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 13.65680433333333
$ws.Range("H17").Value = 40.970413
$ws.Range("I17").Value = 0.154034782011531
$ws.Range("J17").Value = 0.154034782011531
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 144.2727966666667
$ws.Range("N17").Value = 432.81839
$ws.Range("O17").Value = 0.3059077790604178
$ws.Range("P17").Value = 0.3059077790604179
$ws.Range("Q17").Value = 1970.305354699452
$ws.Range("R17").Value = 17732.74819229507
$ws.Range("S17").Value = 0.04712043806320305
$ws.Range("T17").Value = 0.04712043806320307
